# Apply the changes described by the diff:
# 1. Rename the worksheet tab from "alpha1F-HW20.xpc" to "alpha1F"
# 2. Update a handful of floating point values in row 13 (ULP-level corrections)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "alpha1F"

# Update the recalculated values in row 13
$ws.Range("D13").Value = 0.9956681684574961
$ws.Range("E13").Value = 0.9960100990311165
$ws.Range("G13").Value = 0.9986732219721116
$ws.Range("J13").Value = 0.9956681684574961
$ws.Range("K13").Value = 0.9958391337443062
$ws.Range("L13").Value = 0.9970589696431092
$ws.Range("M13").Value = 0.9964040590576269
